$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Changan UNI-S: URL fix (auto_19722 -> auto_19725)
$ws.Range("E9").Value = "https://alyans-auto.ru/auto/auto_19725.html"
$ws.Range("G9").Value = "https://alyans-auto.ru/auto/auto_19725.html"

# Row 11 - Changan UNI-V: URL fix (auto_19626 -> auto_19627)
$ws.Range("E11").Value = "https://alyans-auto.ru/auto/auto_19627.html"
$ws.Range("G11").Value = "https://alyans-auto.ru/auto/auto_19627.html"

# Row 22 - Haval Jolion: URL fix (auto_19911 -> auto_19912)
$ws.Range("E22").Value = "https://alyans-auto.ru/auto/auto_19912.html"
$ws.Range("G22").Value = "https://alyans-auto.ru/auto/auto_19912.html"

# Row 23 - Haval M6: URL fix (auto_20130 -> auto_20131)
$ws.Range("E23").Value = "https://alyans-auto.ru/auto/auto_20131.html"
$ws.Range("G23").Value = "https://alyans-auto.ru/auto/auto_20131.html"

# Row 31 - JAECOO J7: price + URL change (2339900 -> 2639900, auto_20256 -> auto_19844)
$ws.Range("D31").Value = 2639900
$ws.Range("E31").Value = "https://alyans-auto.ru/auto/auto_19844.html"
$ws.Range("F31").Value = 2639900
$ws.Range("G31").Value = "https://alyans-auto.ru/auto/auto_19844.html"

# Row 34 - Lada 4x4 3 dveri: price + URL change (1059000 -> 1055500, auto_19270 -> auto_20482)
$ws.Range("D34").Value = 1055500
$ws.Range("E34").Value = "https://alyans-auto.ru/auto/auto_20482.html"
$ws.Range("F34").Value = 1055500
$ws.Range("G34").Value = "https://alyans-auto.ru/auto/auto_20482.html"
